$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (drop the trailing "...")
$ws.Name = "Basic Clinic Data"

# Move the active selection from D10 to F12 (as last saved in the workbook)
$ws.Range("F12").Select()
